$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114, pushing existing rows 114-133 down to 115-134
$ws.Rows.Item(114).Insert()

# Populate the new row 114 with the new weekly price-report record
$ws.Cells.Item(114, 1).Value = 11
$ws.Cells.Item(114, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(114, 3).Value = "Bíobío"
$ws.Cells.Item(114, 4).Value = 44504
$ws.Cells.Item(114, 5).Value = 8
$ws.Cells.Item(114, 6).Value = 100114001
$ws.Cells.Item(114, 7).Value = "Papa"
$ws.Cells.Item(114, 8).Value = "Patagonia"
$ws.Cells.Item(114, 9).Value = "1a (guarda lavada)"
$ws.Cells.Item(114, 10).Value = 450
$ws.Cells.Item(114, 11).Value = 10000
$ws.Cells.Item(114, 12).Value = 11000
$ws.Cells.Item(114, 13).Value = 10556
$ws.Cells.Item(114, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(114, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(114, 16).Value = 422
$ws.Cells.Item(114, 17).Value = 25
$ws.Cells.Item(114, 18).Value = "Hortaliza"
